# -----------------------------------------------------------------------
# Adds the "Combination" sheet (Citta x Cetasika combination table),
# two workbook-level defined names (CetasikaID / CittaID) used by new
# data-validation list dropdowns, re-points the Citta / Cetasika sheet
# selections, and makes "Combination" the active tab.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- 1. Defined names --------------------------------------------------
$wb.Names.Add('CetasikaID', '=Cetasika!$A$2:$A$53')
$wb.Names.Add('CittaID', '=Citta!$A$2:$A$122')

# Auto-generated data rows for the "Combination" sheet (A2:D93)
$comboData = @(
    @('lobha3','','',''),
    @('lobha3','mana','',''),
    @('lobha4','','',''),
    @('lobha4','mana','',''),
    @('lobha7','','',''),
    @('lobha7','mana','',''),
    @('lobha8','','',''),
    @('lobha8','mana','',''),
    @('dosa1','','',''),
    @('dosa1','issa','',''),
    @('dosa1','macchariya','',''),
    @('dosa1','kukuccha','',''),
    @('dosa2','','',''),
    @('dosa2','issa','',''),
    @('dosa2','macchariya','',''),
    @('dosa2','kukuccha','',''),
    @('dosa2','thina','middha',''),
    @('dosa2','thina','middha','issa'),
    @('dosa2','thina','middha','macchariya'),
    @('dosa2','thina','middha','kukuccha'),
    @('mkus1','','',''),
    @('mkus1','sammavaca','',''),
    @('mkus1','sammakammanta','',''),
    @('mkus1','sammajiva','',''),
    @('mkus1','karuna','',''),
    @('mkus1','mudita','',''),
    @('mkus2','','',''),
    @('mkus2','sammavaca','',''),
    @('mkus2','sammakammanta','',''),
    @('mkus2','sammajiva','',''),
    @('mkus2','karuna','',''),
    @('mkus2','mudita','',''),
    @('mkus3','','',''),
    @('mkus3','sammavaca','',''),
    @('mkus3','sammakammanta','',''),
    @('mkus3','sammajiva','',''),
    @('mkus3','karuna','',''),
    @('mkus3','mudita','',''),
    @('mkus4','','',''),
    @('mkus4','sammavaca','',''),
    @('mkus4','sammakammanta','',''),
    @('mkus4','sammajiva','',''),
    @('mkus4','karuna','',''),
    @('mkus4','mudita','',''),
    @('mkus5','','',''),
    @('mkus5','sammavaca','',''),
    @('mkus5','sammakammanta','',''),
    @('mkus5','sammajiva','',''),
    @('mkus5','karuna','',''),
    @('mkus5','mudita','',''),
    @('mkus6','','',''),
    @('mkus6','sammavaca','',''),
    @('mkus6','sammakammanta','',''),
    @('mkus6','sammajiva','',''),
    @('mkus6','karuna','',''),
    @('mkus6','mudita','',''),
    @('mkus7','','',''),
    @('mkus7','sammavaca','',''),
    @('mkus7','sammakammanta','',''),
    @('mkus7','sammajiva','',''),
    @('mkus7','karuna','',''),
    @('mkus7','mudita','',''),
    @('mkus8','','',''),
    @('mkus8','sammavaca','',''),
    @('mkus8','sammakammanta','',''),
    @('mkus8','sammajiva','',''),
    @('mkus8','karuna','',''),
    @('mkus8','mudita','',''),
    @('mkir1','','',''),
    @('mkir1','karuna','',''),
    @('mkir1','mudita','',''),
    @('mkir2','','',''),
    @('mkir2','karuna','',''),
    @('mkir2','mudita','',''),
    @('mkir3','','',''),
    @('mkir3','karuna','',''),
    @('mkir3','mudita','',''),
    @('mkir4','','',''),
    @('mkir4','karuna','',''),
    @('mkir4','mudita','',''),
    @('mkir5','','',''),
    @('mkir5','karuna','',''),
    @('mkir5','mudita','',''),
    @('mkir6','','',''),
    @('mkir6','karuna','',''),
    @('mkir6','mudita','',''),
    @('mkir7','','',''),
    @('mkir7','karuna','',''),
    @('mkir7','mudita','',''),
    @('mkir8','','',''),
    @('mkir8','karuna','',''),
    @('mkir8','mudita','','')
)


# ---- 2. New "Combination" worksheet, appended after the last sheet ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws.Name = "Combination"

# Header row
$ws.Cells.Item(1, 1).Value = "id"
$ws.Cells.Item(1, 2).Value = "cetasika/0"
$ws.Cells.Item(1, 3).Value = "cetasika/1"
$ws.Cells.Item(1, 4).Value = "cetasika/2"

# Data rows (A2:D93)
$r = 2
foreach ($row in $comboData) {
    if ($row[0] -ne '') { $ws.Cells.Item($r, 1).Value = $row[0] }
    if ($row[1] -ne '') { $ws.Cells.Item($r, 2).Value = $row[1] }
    if ($row[2] -ne '') { $ws.Cells.Item($r, 3).Value = $row[2] }
    if ($row[3] -ne '') { $ws.Cells.Item($r, 4).Value = $row[3] }
    $r++
}

# ---- 3. Data validation (list) dropdowns -------------------------------
$ws.Range("A1:A100").Validation.Add(3, 3, 1, "=CittaID")

$ws.Range("F18:F21").Validation.Add(3, 3, 1, "=CetasikaID")
$ws.Range("B18:D21").Validation.Add(3, 3, 1, "=CetasikaID")
$ws.Range("B2:F17").Validation.Add(3, 3, 1, "=CetasikaID")
$ws.Range("C22:F40").Validation.Add(3, 3, 1, "=CetasikaID")
$ws.Range("B22:B93").Validation.Add(3, 3, 1, "=CetasikaID")

# ---- 4. Re-point selection on "Citta" ----------------------------------
$wsCitta = $wb.Worksheets.Item("Citta")
$wsCitta.Activate()
$wsCitta.Range("A2:A122").Select()

# ---- 5. Re-point selection on "Cetasika" -------------------------------
$wsCetasika = $wb.Worksheets.Item("Cetasika")
$wsCetasika.Activate()
$wsCetasika.Range("A2:A53").Select()
$wsCetasika.Range("A53").Activate()

# ---- 6. Selection on the new "Combination" sheet, make it the active --
#         tab (matches activeTab="6" / tabSelected="1" in the saved file)
$ws.Activate()
$ws.Range("G72").Select()
